# Update countries & provincias Spain
# - Refresh the "last updated" timestamp in A1
# - Swap Hong Kong / Republica de Africa Central ordering (Hong Kong now
#   has more cases than Republica de Africa Central, so it moved above it
#   in the sorted table) and refresh their statistics
# - Refresh the statistics (Casos totales, Nuevos casos, Casos activos,
#   Recuperados, Casos criticos, Muertes) for the other countries whose
#   numbers changed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / timestamp row
$ws.Range("A1").Value = "Datos actualizados a 23 de Agosto de 2020 a las 11:48"

# Row 25 - Filipinas
$ws.Cells.Item(25, 2).Value = 189601
$ws.Cells.Item(25, 3).Value = 2352
$ws.Cells.Item(25, 4).Value = 131367
$ws.Cells.Item(25, 5).Value = 55236
$ws.Cells.Item(25, 7).Value = 32
$ws.Cells.Item(25, 8).Value = 2998

# Row 33 - Israel
$ws.Cells.Item(33, 2).Value = 102150
$ws.Cells.Item(33, 3).Value = 217
$ws.Cells.Item(33, 5).Value = 22022
$ws.Cells.Item(33, 7).Value = 6
$ws.Cells.Item(33, 8).Value = 825

# Row 47 - Polonia
$ws.Cells.Item(47, 2).Value = 61762
$ws.Cells.Item(47, 3).Value = 581
$ws.Cells.Item(47, 4).Value = 42183
$ws.Cells.Item(47, 5).Value = 17624
$ws.Cells.Item(47, 7).Value = 4
$ws.Cells.Item(47, 8).Value = 1955

# Row 49 - Singapur
$ws.Cells.Item(49, 2).Value = 56353
$ws.Cells.Item(49, 3).Value = 87
$ws.Cells.Item(49, 5).Value = 2406

# Row 71 - Austria
$ws.Cells.Item(71, 2).Value = 25253
$ws.Cells.Item(71, 3).Value = 191
$ws.Cells.Item(71, 4).Value = 21558
$ws.Cells.Item(71, 5).Value = 2963

# Row 76 - Estado de Palestina
$ws.Cells.Item(76, 5).Value = 7246
$ws.Cells.Item(76, 7).Value = 2
$ws.Cells.Item(76, 8).Value = 127

# Row 91 - Consejo Danes para los Refugiados
$ws.Cells.Item(91, 2).Value = 9830
$ws.Cells.Item(91, 3).Value = 19
$ws.Cells.Item(91, 4).Value = 8934
$ws.Cells.Item(91, 5).Value = 645

# Row 92 - Malasia
$ws.Cells.Item(92, 2).Value = 9267
$ws.Cells.Item(92, 3).Value = 10
$ws.Cells.Item(92, 4).Value = 8959

# Row 93 - Guinea
$ws.Cells.Item(93, 2).Value = 8967
$ws.Cells.Item(93, 3).Value = 35
$ws.Cells.Item(93, 4).Value = 7708
$ws.Cells.Item(93, 5).Value = 1206

# Row 100 - Finlandia
$ws.Cells.Item(100, 2).Value = 7920
$ws.Cells.Item(100, 3).Value = 14
$ws.Cells.Item(100, 5).Value = 486

# Row 111 now becomes Hong Kong (moved above Republica de Africa Central,
# with refreshed stats)
$ws.Cells.Item(111, 1).Value = "Hong Kong"
$ws.Cells.Item(111, 2).Value = 4683
$ws.Cells.Item(111, 3).Value = 25
$ws.Cells.Item(111, 4).Value = 4018
$ws.Cells.Item(111, 5).Value = 588
$ws.Cells.Item(111, 7).Value = 1
$ws.Cells.Item(111, 8).Value = 77

# Row 112 now becomes Republica de Africa Central (its stats carry over
# unchanged from the old row 111)
$ws.Cells.Item(112, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(112, 2).Value = 4679
$ws.Cells.Item(112, 4).Value = 1755
$ws.Cells.Item(112, 5).Value = 2863
$ws.Cells.Item(112, 8).Value = 61

# Row 125 - Sri Lanka
$ws.Cells.Item(125, 4).Value = 2805
$ws.Cells.Item(125, 5).Value = 130
